$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.199.07"
$ws.Range("D3").Value = "2.520.30"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +4.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0820"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "2.915.01"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "2.526.71"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "48.075.29"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.87%  "
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.21%  "
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("E30").Value = "  +4.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "2.003.62"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.59%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.82%  "
